$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Taxonsorteringsordning changes
$ws.Range("B2").Value = 95707

# Row 3 gets old row 6's identifying data (A/B/Q/R) - species stays "Lunglav"
$ws.Range("A3").Value = 112181635
$ws.Range("B3").Value = 78713
$ws.Range("Q3").Value = 492738
$ws.Range("R3").Value = 6845111

# Row 4 gets old row 3's identifying data + species switches to "Lunglav"
$ws.Range("A4").Value = 112181727
$ws.Range("B4").Value = 78713
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("Q4").Value = 493005
$ws.Range("R4").Value = 6845384
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-07-04"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-07-04"

# Row 5 gets old row 4's identifying data + species switches to "Kolflarnlav"
$ws.Range("A5").Value = 112182531
$ws.Range("B5").Value = 77402
$ws.Range("E5").Value = 6446
$ws.Range("F5").Value = "Kolflarnlav"
$ws.Range("G5").Value = "Carbonicola anthracophila"
$ws.Range("H5").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q5").Value = 492845
$ws.Range("R5").Value = 6845301

# Row 6 gets old row 5's identifying data - species stays "Lunglav"
$ws.Range("A6").Value = 112183278
$ws.Range("B6").Value = 78713
$ws.Range("Q6").Value = 492536
$ws.Range("R6").Value = 6845328
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-07-03"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-07-03"
